$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-PlainValue($cellRef, $val) {
    $ws.Range($cellRef).Value = $val
}

Set-TextValue "D2" "29.551.10"
Set-PlainValue "E2" "  +0.65%  "

Set-TextValue "D3" "1.850.36"
Set-PlainValue "E3" "  +0.15%  "

Set-TextValue "D4" "0.9975"
Set-PlainValue "E4" "  -0.24%  "

Set-TextValue "D5" "240.58"
Set-PlainValue "E5" "  +0.18%  "

Set-TextValue "D6" "0.6307"
Set-PlainValue "E6" "  +0.46%  "

Set-TextValue "D7" "0.9991"
Set-PlainValue "E7" "  -0.11%  "

Set-TextValue "D8" "0.07495"
Set-PlainValue "E8" "  -1.12%  "

Set-TextValue "D9" "0.2917"
Set-PlainValue "E9" "  +0.27%  "

Set-TextValue "D10" "24.70"
Set-PlainValue "E10" "  +0.71%  "

Set-TextValue "D11" "0.07743"
Set-PlainValue "E11" "  -0.09%  "

Set-TextValue "D12" "1.848.34"
Set-PlainValue "E12" "  +0.05%  "

Set-TextValue "D13" "5.021"
Set-PlainValue "E13" "  +0.14%  "

Set-TextValue "D14" "0.6831"
Set-PlainValue "E14" "  +0.65%  "

Set-TextValue "D15" "0.00001047"
Set-PlainValue "E15" "  +0.85%  "

Set-TextValue "D16" "82.31"
Set-PlainValue "E16" "  -0.86%  "

Set-TextValue "D17" "6.266"
Set-PlainValue "E17" "  +2.55%  "

Set-TextValue "D18" "29.529.80"
Set-PlainValue "E18" "  +0.56%  "

Set-TextValue "D19" "230.14"
Set-PlainValue "E19" "  +0.68%  "

Set-TextValue "D20" "12.42"
Set-PlainValue "E20" "  +0.85%  "

Set-TextValue "D21" "0.9991"
Set-PlainValue "E21" "  -0.10%  "

Set-TextValue "D22" "7.564"
Set-PlainValue "E22" "  +1.88%  "

Set-TextValue "D23" "0.9989"
Set-PlainValue "E23" "  -0.18%  "

Set-TextValue "D24" "159.57"
Set-PlainValue "E24" "  +0.52%  "

Set-TextValue "D25" "8.536"
Set-PlainValue "E25" "  +1.25%  "

Set-TextValue "D26" "0.1373"
Set-PlainValue "E26" "  -1.02%  "

Set-TextValue "D27" "17.55"
Set-PlainValue "E27" "  -0.56%  "

Set-TextValue "D28" "0.06541"
Set-PlainValue "E28" "  +15.97%  "

Set-TextValue "D29" "1.428"
Set-PlainValue "E29" "  -0.22%  "

Set-TextValue "D30" "1.489"
Set-PlainValue "E30" "  +1.33%  "

Set-TextValue "D31" "4.108"
Set-PlainValue "E31" "  +0.01%  "

Set-TextValue "D32" "4.108"
Set-PlainValue "E32" "  +1.86%  "

Set-TextValue "D33" "1.839"
Set-PlainValue "E33" "  +0.95%  "

Set-TextValue "D34" "1.149"
Set-PlainValue "E34" "  -0.54%  "

Set-TextValue "D35" "0.6997"
Set-PlainValue "E35" "  +0.56%  "

Set-PlainValue "E36" "  -0.19%  "

Set-TextValue "D37" "0.01868"
Set-PlainValue "E37" "  +2.36%  "

Set-TextValue "D38" "1.266.26"
Set-PlainValue "E38" "  +2.75%  "

Set-PlainValue "E39" "  +4.20%  "

Set-TextValue "D40" "6.846"
Set-PlainValue "E40" "  +7.34%  "

Set-TextValue "D41" "0.9401"
Set-PlainValue "E41" "  +4.58%  "

Set-TextValue "D42" "0.9999"
Set-PlainValue "E42" "  +0.03%  "

Set-TextValue "D43" "2.027.56"
Set-PlainValue "E43" "  +1.10%  "

Set-TextValue "D44" "101.34"
Set-PlainValue "E44" "  +0.19%  "

Set-TextValue "D45" "66.39"
Set-PlainValue "E45" "  +1.70%  "

Set-TextValue "D46" "1.744"
Set-PlainValue "E46" "  +4.15%  "

Set-PlainValue "B47" "Aptos"
Set-PlainValue "C47" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.117"
Set-PlainValue "E47" "  -0.01%  "

Set-PlainValue "B48" "Algorand"
Set-PlainValue "C48" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D48" "0.1169"
Set-PlainValue "E48" "  +1.89%  "

Set-PlainValue "B49" "BabyDogeCoin"
Set-PlainValue "C49" "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.00000000116"
Set-PlainValue "E49" "  +1.29%  "

Set-TextValue "D50" "9.001"
Set-PlainValue "E50" "  +0.09%  "

Set-TextValue "D51" "0.3970"
Set-PlainValue "E51" "  -0.52%  "
